$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo "выполненых" -> "выполненных" ("completed" with correct spelling)
# throughout every cell on the sheet (this text appears many times inside the
# multi-line log cells as well as in several standalone short cells).
$null = $ws.Cells.Replace("выполненых", "выполненных", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)

# The replace above lengthens some of the multi-line log cells, which makes
# the engine auto-expand the row height of the affected rows. Re-run AutoFit
# on just those rows so the row heights collapse back to their original,
# unset (default) state instead of being left with an explicit custom height.
foreach ($r in $ws.UsedRange.Rows) {
    $r.AutoFit()
}
